# 07/10/22 Added Movable Chinese Title + Fixed Hymnal "Hymns" Positioning
#
# Reposition four text boxes on slide 1. The PowerPoint COM object model
# expresses Shape.Top/.Left in points, while the underlying OOXML stores
# EMU (914400 EMU = 1 inch = 72 points, i.e. 12700 EMU per point). A tiny
# (+0.5 EMU) nudge is added before converting to points so that internal
# float-precision truncation when writing back to EMU cannot round the
# result down to one EMU less than intended.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "English" title textbox -> move up (make room for the new movable Chinese title)
$s.Shapes.Item(1).Top = (66328 + 0.5) / 914400 * 72

# "Hymns / 詩" textbox -> nudge down
$s.Shapes.Item(3).Top = (2406367 + 0.5) / 914400 * 72

# "Hymn No." textbox -> nudge down
$s.Shapes.Item(4).Top = (3414480 + 0.5) / 914400 * 72

# "Bible Verse / 經文" textbox -> nudge down
$s.Shapes.Item(5).Top = (2561349 + 0.5) / 914400 * 72
